$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A28").Value = "Springfield"
$ws.Range("B28").Value = "Springs"
$ws.Range("C28").Value = "free conference"
$ws.Range("D28").Value = "Turts#3627"
$ws.Range("E28").Value = "Hunter Scott"
$ws.Range("F28").Value = "air raid"
$ws.Range("G28").Value = "4-4"
$ws.Range("H28").Value = "0-0"
